$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first row ("Требуется КУРЬЕР"), shifting remaining rows up
$ws.Rows("1").Delete()

# Add two new product rows at the bottom
$ws.Range("A3").Value = "Оффлайн ТВ 2 плитки"
$ws.Range("B3").Value = 2300
$ws.Range("C3").Value = 167

$ws.Range("A4").Value = "Оффлайн ТВ 3 плитки"
$ws.Range("B4").Value = 2800
$ws.Range("C4").Value = 170

# Restore the selection on the row after the data
$ws.Range("A4").Select()
